$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177, pushing the existing rows 177:210 down to 178:211.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row with the new weekly price observation.
$ws.Cells.Item(177, 1).Value = 10
$ws.Cells.Item(177, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(177, 3).Value = "La Araucanía"
$ws.Cells.Item(177, 4).Value = 45173
$ws.Cells.Item(177, 5).Value = 9
$ws.Cells.Item(177, 6).Value = 100112035
$ws.Cells.Item(177, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 80
$ws.Cells.Item(177, 11).Value = 25000
$ws.Cells.Item(177, 12).Value = 25000
$ws.Cells.Item(177, 13).Value = 25000
$ws.Cells.Item(177, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(177, 15).Value = "Región Metropolitana"
$ws.Cells.Item(177, 16).Value = 1667
$ws.Cells.Item(177, 17).Value = 15
$ws.Cells.Item(177, 18).Value = "Hortaliza"
